# Add the remaining "4-character word" rows beneath the existing
# header row (単語 / 意味) on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$words = @("たいよう", "せいたい", "かんたい", "かいたい", "こうせい")

$row = 2
foreach ($word in $words) {
    $ws.Cells.Item($row, 1).Value = $word
    $row = $row + 1
}
